$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 does not exist yet; copy the formatting of row 18 down to row 19 first so the
# new row picks up the same style (s="4") as the rest of the data rows.
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the old row 18 (K_URBAN) down to row 19
$ws.Range("A19").Value = $ws.Range("A18").Value2
$ws.Range("B19").Value = $ws.Range("B18").Value2
$ws.Range("C19").Value = $ws.Range("C18").Value2

# Row 18 becomes the former K_TYPEAREA row (currently still in row 17)
$ws.Range("A18").Value = $ws.Range("A17").Value2
$ws.Range("B18").Value = $ws.Range("B17").Value2
$ws.Range("C18").Value = $ws.Range("C17").Value2

# Row 17 becomes the new K_TARIF entry
$ws.Range("A17").Value = "K_TARIF"
$ws.Range("B17").Value = "Art der Tarifvertr" + [char]0x00E4 + "ge"
$ws.Range("C17").Value = "XXXArt der Tarifvertr" + [char]0x00E4 + "ge"
